$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("orders")
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "products"
$ws1.Range("A1:AA2").Copy($ws2.Range("A1"))
$r = $ws2.Range("U1:X1").EntireColumn
$r.Insert()
